$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF" (same style as existing headers) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from H1 (bold, bordered, centered) onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-30 for columns I (I0) and J (IF) ---
$values = @(
    @(5, 6),
    @(6, 8),
    @(8, 8),
    @(7, 7),
    @(10, 10),
    @(9, 9),
    @(2, 4),
    @(4, 7),
    @(9, 10),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(6, 7),
    @(9, 9),
    @(6, 8),
    @(6, 8),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(5, 6),
    @(5, 7),
    @(6, 7),
    @(6, 7),
    @(6, 8),
    @(1, 3),
    @(5, 7),
    @(3, 4)
)

for ($idx = 0; $idx -lt $values.Count; $idx++) {
    $row = 2 + $idx
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}

Write-Host "Added I0/IF columns"
